$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1598
$ws.Range("I6").Value = 409.5
$ws.Range("K6").Value = 1228.5
$ws.Range("M6").Value = -1116.5
$ws.Range("H9").Value = 209.91667
$ws.Range("I9").Value = 232
$ws.Range("J9").Value = 99.5
$ws.Range("K9").Value = 232
$ws.Range("L9").Value = 99.5
$ws.Range("M9").Value = -63
$ws.Range("N9").Value = -437.5
$ws.Range("H12").Value = 38710.152
$ws.Range("I12").Value = 316
$ws.Range("J12").Value = 62706.5
$ws.Range("K12").Value = 316
$ws.Range("L12").Value = 62706.5
$ws.Range("M12").Value = -146
$ws.Range("N12").Value = -63046.5
$ws.Range("H21").Value = 12500
$ws.Range("I21").Value = 12500
$ws.Range("J21").Value = 12500
$ws.Range("K21").Value = 12500
$ws.Range("L21").Value = 12500
$ws.Range("M21").Value = -12032
$ws.Range("N21").Value = -13436
$ws.Range("H23").Value = 12500
$ws.Range("I23").Value = 12500
$ws.Range("J23").Value = 12500
$ws.Range("K23").Value = 12500
$ws.Range("L23").Value = 12500
$ws.Range("M23").Value = -12266
$ws.Range("N23").Value = -12968
$ws.Range("H29").Value = 812.2727
$ws.Range("J29").Value = 812.2727
$ws.Range("L29").Value = 2436.8181
$ws.Range("N29").Value = -2998.8181
$ws.Range("H38").Value = 2016549.5
$ws.Range("I38").Value = 5376496.5
$ws.Range("J38").Value = 581.5
$ws.Range("K38").Value = 16129489.5
$ws.Range("L38").Value = 1744.5
$ws.Range("M38").Value = -16129117.5
$ws.Range("N38").Value = -2488.5
$ws.Range("H58").Value = 722626.1
$ws.Range("I58").Value = 1165734.5
$ws.Range("J58").Value = 2575
$ws.Range("K58").Value = 3497203.5
$ws.Range("L58").Value = 7725
$ws.Range("M58").Value = -3497053.5
$ws.Range("N58").Value = -8025
$ws.Range("H62").Value = 2221.1667
$ws.Range("I62").Value = 2185
$ws.Range("J62").Value = 2329.6667
$ws.Range("K62").Value = 2185
$ws.Range("L62").Value = 2329.6667
$ws.Range("M62").Value = -1561
$ws.Range("N62").Value = -3577.6667
$ws.Range("H65").Value = 2221.1667
$ws.Range("I65").Value = 2185
$ws.Range("J65").Value = 2329.6667
$ws.Range("K65").Value = 10925
$ws.Range("L65").Value = 11648.3335
$ws.Range("M65").Value = -7805
$ws.Range("N65").Value = -17888.3335
$ws.Range("H87").Value = 31222.111
$ws.Range("J87").Value = 31222.111
$ws.Range("L87").Value = 31222.111
$ws.Range("N87").Value = -33718.111
$ws.Range("H90").Value = 31222.111
$ws.Range("J90").Value = 31222.111
$ws.Range("L90").Value = 93666.333
$ws.Range("N90").Value = -106146.333
$ws.Range("H132").Value = 4633729.5
$ws.Range("I132").Value = 4811922.5
$ws.Range("J132").Value = 703
$ws.Range("K132").Value = 14435767.5
$ws.Range("L132").Value = 2109
$ws.Range("M132").Value = -14433237.5
$ws.Range("N132").Value = -7169
$ws.Range("H135").Value = 446.64517
$ws.Range("I135").Value = 461.0345
$ws.Range("J135").Value = 238
$ws.Range("K135").Value = 4149.3105
$ws.Range("L135").Value = 2142
$ws.Range("M135").Value = -1614.3105
$ws.Range("N135").Value = -7212
$ws.Range("H137").Value = 1383.8431
$ws.Range("I137").Value = 1580.6666
$ws.Range("J137").Value = 1323.2821
$ws.Range("K137").Value = 4741.9998
$ws.Range("L137").Value = 3969.8463
$ws.Range("M137").Value = -2191.9998
$ws.Range("N137").Value = -9069.846299999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 25785.098
$ws.Range("I2").Value = 1172.0358
$ws.Range("J2").Value = 78797.84
$ws.Range("K2").Value = 1172.0358
$ws.Range("L2").Value = 78797.84
$ws.Range("M2").Value = -1059.0358
$ws.Range("N2").Value = -79023.84
$ws.Range("H32").Value = 22398.186
$ws.Range("I32").Value = 4002
$ws.Range("K32").Value = 4002
$ws.Range("M32").Value = -3715
$ws.Range("H110").Value = 19270894
$ws.Range("I110").Value = 27834498
$ws.Range("J110").Value = 2789
$ws.Range("K110").Value = 27834498
$ws.Range("L110").Value = 2789
$ws.Range("M110").Value = -27832453
$ws.Range("N110").Value = -6879
$ws.Range("H116").Value = 25785.098
$ws.Range("I116").Value = 1172.0358
$ws.Range("J116").Value = 78797.84
$ws.Range("K116").Value = 1172.0358
$ws.Range("L116").Value = 78797.84
$ws.Range("M116").Value = 1121.9642
$ws.Range("N116").Value = -83385.84
$ws.Range("H122").Value = 4051
$ws.Range("I122").Value = 4690.2
$ws.Range("J122").Value = 2453
$ws.Range("K122").Value = 14070.6
$ws.Range("L122").Value = 7359
$ws.Range("M122").Value = -11620.6
$ws.Range("N122").Value = -12259

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 25785.098
$ws.Range("I3").Value = 1172.0358
$ws.Range("J3").Value = 78797.84
$ws.Range("K3").Value = 1172.0358
$ws.Range("L3").Value = 78797.84
$ws.Range("M3").Value = -1058.0358
$ws.Range("N3").Value = -79025.84
$ws.Range("H43").Value = 318500.25
$ws.Range("J43").Value = 318500.25
$ws.Range("L43").Value = 318500.25
$ws.Range("N43").Value = -318862.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16673.383
$ws.Range("I31").Value = 33290.13
$ws.Range("J31").Value = 2751.2432
$ws.Range("K31").Value = 33290.13
$ws.Range("L31").Value = 2751.2432
$ws.Range("M31").Value = -32995.13
$ws.Range("N31").Value = -3341.2432
$ws.Range("H34").Value = 16673.383
$ws.Range("I34").Value = 33290.13
$ws.Range("J34").Value = 2751.2432
$ws.Range("K34").Value = 33290.13
$ws.Range("L34").Value = 2751.2432
$ws.Range("M34").Value = -33088.13
$ws.Range("N34").Value = -3155.2432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2166.6667
$ws.Range("H34").Value = 2181.125
$ws.Range("I34").Value = 150
$ws.Range("J34").Value = 2471.2856
$ws.Range("K34").Value = 450
$ws.Range("L34").Value = 7413.8568
$ws.Range("M34").Value = -366
$ws.Range("N34").Value = -7581.8568
$ws.Range("H39").Value = 3466.6667
$ws.Range("J39").Value = 3466.6667
$ws.Range("L39").Value = 10400.0001
$ws.Range("N39").Value = -10988.0001
$ws.Range("H55").Value = 8287.833000000001
$ws.Range("J55").Value = 8287.833000000001
$ws.Range("L55").Value = 24863.499
$ws.Range("N55").Value = -25217.499
$ws.Range("H113").Value = 559.7568
$ws.Range("I113").Value = 517.25
$ws.Range("J113").Value = 592.1429000000001
$ws.Range("K113").Value = 1551.75
$ws.Range("L113").Value = 1776.4287
$ws.Range("M113").Value = 618.25
$ws.Range("N113").Value = -6116.4287
$ws.Range("H131").Value = 1384.9259
$ws.Range("I131").Value = 1127.1428
$ws.Range("J131").Value = 1423.3191
$ws.Range("K131").Value = 3381.4284
$ws.Range("L131").Value = 4269.9573
$ws.Range("M131").Value = 1658.5716
$ws.Range("N131").Value = -14349.9573

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2306.3333
$ws.Range("I43").Value = 1950
$ws.Range("J43").Value = 3019
$ws.Range("K43").Value = 1950
$ws.Range("L43").Value = 3019
$ws.Range("M43").Value = -1799
$ws.Range("N43").Value = -3321
$ws.Range("H113").Value = 1529.3077
$ws.Range("I113").Value = 1220.1666
$ws.Range("J113").Value = 1794.2858
$ws.Range("K113").Value = 1220.1666
$ws.Range("L113").Value = 1794.2858
$ws.Range("M113").Value = 949.8334
$ws.Range("N113").Value = -6134.2858
$ws.Range("H131").Value = 35325
$ws.Range("I131").Value = 35325
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 35325
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -30285
$ws.Range("N131").Value = $null
$ws.Range("H132").Value = 2947.7273
$ws.Range("I132").Value = 2433.5625
$ws.Range("K132").Value = 7300.6875
$ws.Range("M132").Value = -4770.6875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 202396
$ws.Range("I40").Value = 501000
$ws.Range("K40").Value = 501000
$ws.Range("M40").Value = -500864
$ws.Range("H61").Value = 2168.8572
$ws.Range("I61").Value = 2066.4
$ws.Range("J61").Value = 2425
$ws.Range("K61").Value = 2066.4
$ws.Range("L61").Value = 2425
$ws.Range("M61").Value = -1864.4
$ws.Range("N61").Value = -2829
$ws.Range("H104").Value = 17222
$ws.Range("J104").Value = 17222
$ws.Range("L104").Value = 17222
$ws.Range("N104").Value = -24210
$ws.Range("H113").Value = 2168.8572
$ws.Range("I113").Value = 2066.4
$ws.Range("J113").Value = 2425
$ws.Range("K113").Value = 2066.4
$ws.Range("L113").Value = 2425
$ws.Range("M113").Value = 103.5999999999999
$ws.Range("N113").Value = -6765
$ws.Range("H122").Value = 3498.5
$ws.Range("I122").Value = 3498.5
$ws.Range("J122").Value = 3498.5
$ws.Range("K122").Value = 10495.5
$ws.Range("L122").Value = 10495.5
$ws.Range("M122").Value = -8045.5
$ws.Range("N122").Value = -15395.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H102").Value = 46158.5
$ws.Range("J102").Value = 46158.5
$ws.Range("L102").Value = 46158.5
$ws.Range("N102").Value = -52648.5
$ws.Range("H103").Value = 15481.2
$ws.Range("J103").Value = 15481.2
$ws.Range("L103").Value = 15481.2
$ws.Range("N103").Value = -17825.2
$ws.Range("H104").Value = 23266.666
$ws.Range("J104").Value = 23266.666
$ws.Range("L104").Value = 23266.666
$ws.Range("N104").Value = -30254.666
$ws.Range("H106").Value = 18550
$ws.Range("J106").Value = 18550
$ws.Range("L106").Value = 18550
$ws.Range("N106").Value = -21074
$ws.Range("H131").Value = 52500
$ws.Range("J131").Value = 52500
$ws.Range("L131").Value = 52500
$ws.Range("N131").Value = -62580
